# Daily attendance processing - reorder "Recorded By" (column G) entries.
# For each data row, the comma-separated list of recorders is reversed,
# except for the special case "System, admin@admin.com" which is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    if ($val -eq "System, admin@admin.com") {
        continue
    }

    if ($val -like "*,*") {
        $parts = $val -split ",\s*"
        $reversedParts = $parts[($parts.Count - 1)..0]
        $newVal = [string]::Join(", ", $reversedParts)
        $cell.Value2 = $newVal
    }
}
